$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (I0, IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold, border, centered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-49
$data = @(
    @(7,8),
    @(7,9),
    @(5,7),
    @(7,7),
    @(8,8),
    @(10,10),
    @(8,8),
    @(1,6),
    @(1,2),
    @(1,2),
    @(1,6),
    @(1,6),
    @(1,5),
    @(6,7),
    @(1,3),
    @(9,9),
    @(1,3),
    @(1,3),
    @(7,7),
    @(8,8),
    @(7,9),
    @(4,6),
    @(7,9),
    @(4,6),
    @(8,8),
    @(5,6),
    @(14,14),
    @(8,9),
    @(5,5),
    @(7,9),
    @(9,9),
    @(7,8),
    @(5,6),
    @(8,9),
    @(9,9),
    @(4,6),
    @(6,9),
    @(7,8),
    @(6,9),
    @(9,9),
    @(2,6),
    @(7,7),
    @(8,9),
    @(6,7),
    @(1,4),
    @(4,6),
    @(3,4),
    @(4,5)
)

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
